$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank columns before column D, shifting existing D:K data to F:M
$ws.Range("D5:E102").Insert(-4161, 0)

# Copy number/date formatting from the (shifted) original column D -- now column F -- into the two new columns
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("F7:F35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the newest two quarters of data
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 23381000
$ws.Range("E8").Value2 = 23257000
$ws.Range("D9").Value2 = "NA"
$ws.Range("E9").Value2 = "NA"
$ws.Range("D10").Value2 = "NA"
$ws.Range("E10").Value2 = "NA"
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 2000
$ws.Range("E14").Value2 = 5000
$ws.Range("D15").Value2 = 93000
$ws.Range("E15").Value2 = 91000
$ws.Range("D17").Value2 = 22649000
$ws.Range("E17").Value2 = 21827000
$ws.Range("D18").Value2 = 732000
$ws.Range("E18").Value2 = 1430000
$ws.Range("D20").Value2 = 0
$ws.Range("E20").Value2 = 0
$ws.Range("D21").Value2 = 1020000
$ws.Range("E21").Value2 = 1718300
$ws.Range("D22").Value2 = 189000
$ws.Range("E22").Value2 = 188000
$ws.Range("D23").Value2 = 543000
$ws.Range("E23").Value2 = 1242000
$ws.Range("D24").Value2 = 146000
$ws.Range("E24").Value2 = 282000
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 397000
$ws.Range("E26").Value2 = 960000
$ws.Range("D27").Value2 = 397000
$ws.Range("E27").Value2 = 960000
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 28000
$ws.Range("E29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("E32").Value2 = 0
$ws.Range("D33").Value2 = 425000
$ws.Range("E33").Value2 = 960000
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 425000
$ws.Range("E35").Value2 = 960000
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 3934000
$ws.Range("E41").Value2 = 4260000
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 6743000
$ws.Range("E43").Value2 = 6943000
$ws.Range("D44").Value2 = 0
$ws.Range("E44").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("E45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 0
$ws.Range("D47").Value2 = 25020000
$ws.Range("E47").Value2 = 26228000
$ws.Range("D48").Value2 = 2735000
$ws.Range("E48").Value2 = 2592000
$ws.Range("D49").Value2 = 29511000
$ws.Range("E49").Value2 = 29569000
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 604000
$ws.Range("E52").Value2 = 741000
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 71571000
$ws.Range("E54").Value2 = 74445000
$ws.Range("D57").Value2 = 4959000
$ws.Range("E57").Value2 = 6286000
$ws.Range("D58").Value2 = 1994000
$ws.Range("E58").Value2 = 2119000
$ws.Range("D59").Value2 = 2484000
$ws.Range("E59").Value2 = 2545000
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 0
$ws.Range("D61").Value2 = 17217000
$ws.Range("E61").Value2 = 17300000
$ws.Range("D62").Value2 = 1960000
$ws.Range("E62").Value2 = 2063000
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 43030000
$ws.Range("E66").Value2 = 45351000
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 19988000
$ws.Range("E72").Value2 = 20182000
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 28541000
$ws.Range("E76").Value2 = 29094000
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 425000
$ws.Range("E81").Value2 = 960000
$ws.Range("D83").Value2 = 288000
$ws.Range("E83").Value2 = 288300
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 463000
$ws.Range("E89").Value2 = 606900
$ws.Range("D91").Value2 = -320000
$ws.Range("E91").Value2 = -355500
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = 773000
$ws.Range("E94").Value2 = -687800
$ws.Range("D96").Value2 = -193000
$ws.Range("E96").Value2 = -194700
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -1561000
$ws.Range("E100").Value2 = -340600
$ws.Range("D101").Value2 = -1000
$ws.Range("E101").Value2 = -600
$ws.Range("D102").Value2 = -326000
$ws.Range("E102").Value2 = -422100
